$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Install" cell display text to "Link" (hyperlink itself is unchanged) ---
$ws.Range("C2").Value2 = "Link"

# --- Rework the Terabox hyperlink cell (D2): shorten its displayed text to "Terabox" ---
# Remove the existing hyperlink on D2 (it currently displays the full Terabox URL)
$existingStyle = $ws.Range("D2").Style
$toDelete = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$2') {
        $toDelete += $h
    }
}
foreach ($h in $toDelete) {
    $h.Delete()
}

# Set the new short label text
$ws.Range("D2").Value2 = "Terabox"

# Re-add the hyperlink on D2, pointing at the same Terabox share URL as before
$ws.Hyperlinks.Add($ws.Range("D2"), "https://1024terabox.com/s/1H0j6ZcIszFVzd4882qy-4A") | Out-Null

# Restore the original (Hyperlink) cell style so D2 matches the other link cells
$ws.Range("D2").Style = $existingStyle

# --- Update the view: scroll so column C is left-most visible, move the selection to D10 ---
$ws.Range("D10").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
